$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2 through 97) forward by 7 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 7
}

# Update the Actual Production values in column B (rows 29 through 49)
$newValues = @{
    29 = 12
    30 = 43
    31 = 88
    32 = 147
    33 = 224
    34 = 299
    35 = 391
    36 = 471
    37 = 563
    38 = 648
    39 = 757
    40 = 836
    41 = 911
    42 = 974
    43 = 996
    44 = 1107
    45 = 1097
    46 = 1079
    47 = 1083
    48 = 1106
    49 = 1083
}

foreach ($r in $newValues.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $newValues[$r]
}
